$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for columns BI (median_temp) and BJ (mad_temp),
# matching the bold/border header style used by the rest of row 1.
$ws.Range("BI1").Value = "median_temp"
$ws.Range("BJ1").Value = "mad_temp"
$ws.Range("BH1").Copy()
$ws.Range("BI1:BJ1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New data values for rows 2-39 (median_temp / mad_temp per row)
$ws.Range("BI2").Value = 32.01123046875
$ws.Range("BJ2").Value = 0.0015869140625
$ws.Range("BI3").Value = 32.03125
$ws.Range("BJ3").Value = 0.001035690307617188
$ws.Range("BI4").Value = 32.0458984375
$ws.Range("BJ4").Value = 0.001071929931640625
$ws.Range("BI5").Value = 32.0458984375
$ws.Range("BJ5").Value = 0.00136566162109375
$ws.Range("BI6").Value = 32.072265625
$ws.Range("BJ6").Value = 0.001230239868164062
$ws.Range("BI7").Value = 32.08984375
$ws.Range("BJ7").Value = 0.00092315673828125
$ws.Range("BI8").Value = 32.1611328125
$ws.Range("BJ8").Value = 0.001279830932617188
$ws.Range("BI9").Value = 32.1875
$ws.Range("BJ9").Value = 0.00160980224609375
$ws.Range("BI10").Value = 32.2158203125
$ws.Range("BJ10").Value = 0.0006694793701171875
$ws.Range("BI11").Value = 32.5439453125
$ws.Range("BJ11").Value = 0.001209259033203125
$ws.Range("BI12").Value = 32.5576171875
$ws.Range("BJ12").Value = 0.0009918212890625
$ws.Range("BI13").Value = 32.5556640625
$ws.Range("BJ13").Value = 0.0006256103515625
$ws.Range("BI14").Value = 32.6796875
$ws.Range("BJ14").Value = 0.002063751220703125
$ws.Range("BI15").Value = 32.7216796875
$ws.Range("BJ15").Value = 0.00153350830078125
$ws.Range("BI16").Value = 32.7509765625
$ws.Range("BJ16").Value = 0.001036644857668042
$ws.Range("BI17").Value = 32.890625
$ws.Range("BJ17").Value = 0.00147247314453125
$ws.Range("BI18").Value = 32.91796875
$ws.Range("BJ18").Value = 0.001129907024793114
$ws.Range("BI19").Value = 32.93359375
$ws.Range("BJ19").Value = 0.0008640289306640625
$ws.Range("BI20").Value = 32.92578125
$ws.Range("BJ20").Value = 0.000962857910819446
$ws.Range("BI21").Value = 32.9345703125
$ws.Range("BJ21").Value = 0.0009398695676772962
$ws.Range("BI22").Value = 32.94140625
$ws.Range("BJ22").Value = 0.0008744412064287948
$ws.Range("BI23").Value = 32.94921875
$ws.Range("BJ23").Value = 0.0008266961860571841
$ws.Range("BI24").Value = 32.5029296875
$ws.Range("BJ24").Value = 0.01142811849635794
$ws.Range("BI25").Value = 32.34765625
$ws.Range("BJ25").Value = 0.004897039216961543
$ws.Range("BI26").Value = 32.287109375
$ws.Range("BJ26").Value = 0.002666170634920748
$ws.Range("BI27").Value = 32.2470703125
$ws.Range("BJ27").Value = 0.002535912298387097
$ws.Range("BI28").Value = 32.0498046875
$ws.Range("BJ28").Value = 0.003091832797933997
$ws.Range("BI29").Value = 32.0908203125
$ws.Range("BJ29").Value = 0.001450108936004362
$ws.Range("BI30").Value = 32.1025390625
$ws.Range("BJ30").Value = 0.001015684051399009
$ws.Range("BI31").Value = 32.1005859375
$ws.Range("BJ31").Value = 0.0009464676030118749
$ws.Range("BI32").Value = 32.115234375
$ws.Range("BJ32").Value = 0.001104043131059852
$ws.Range("BI33").Value = 32.1162109375
$ws.Range("BJ33").Value = 0.0008812056353878903
$ws.Range("BI34").Value = 32.119140625
$ws.Range("BJ34").Value = 0.001055875392982898
$ws.Range("BI35").Value = 32.119140625
$ws.Range("BJ35").Value = 0.0009832145490892535
$ws.Range("BI36").Value = 32.1083984375
$ws.Range("BJ36").Value = 0.003790113940693439
$ws.Range("BI37").Value = 32.0830078125
$ws.Range("BJ37").Value = 0.002440547744951888
$ws.Range("BI38").Value = 32.06640625
$ws.Range("BJ38").Value = 0.001983198691992361
$ws.Range("BI39").Value = 32.0478515625
$ws.Range("BJ39").Value = 0.002123311327117022
